$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.306451612903226
$ws.Range("C2").Value = 0.359313077939234
$ws.Range("D2").Value = 0.756906077348066
$ws.Range("E2").Value = 0.277950310559006
$ws.Range("F2").Value = 0.0408626560726447

$ws.Range("B3").Value = 0.52258064516129
$ws.Range("C3").Value = 0.696169088507266
$ws.Range("D3").Value = 0.869244935543278
$ws.Range("E3").Value = 0.631987577639752
$ws.Range("F3").Value = 0.348467650397276

$ws.Range("B4").Value = 0.403225806451613
$ws.Range("C4").Value = 0.578599735799207
$ws.Range("D4").Value = 0.69060773480663
$ws.Range("E4").Value = 0.389751552795031
$ws.Range("F4").Value = 0.149829738933031
